# Daily attendance processing - 2025-09-30 20:20:26
# Applies the attendance-report refresh: reshuffled "Recorded By" name
# lists, updated dates/times for several sessions, refreshed attendance
# statistics, and the PHYSIOLOGY A3 session #1 row (91) that moved from
# "Pending" to "Recorded" (including its own style/highlight change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2 (ANATOMY A1 #1) - reorder recorder emails
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 4 (ANATOMY A1 #3) - date shifted a day
$ws.Range("E4").Value = "15/10/2025"

# Class statistics block (K/L columns)
$ws.Range("L6").Value = 19
$ws.Range("L8").Value = 167

# Row 9 (HISTOLOGY A1 #1) - reorder recorder emails
$ws.Range("G9").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# Coverage / average attendance percentages (stored as literal text).
# A leading apostrophe keeps Excel from re-parsing "nn.n%" into a real
# percentage number; the donor-cell PasteSpecial pass below then
# restores each cell's original style index (they must stay s="5").
$ws.Range("L9").Value = "'10.2%"
$ws.Range("L10").Value = "'45.1%"

# Row 13 (MICROBIOLOGY A1 #2) - time changed
$ws.Range("F13").Value = "12:00:00"

# Group statistics row 17 (Year 3 / A3)
$ws.Range("O17").Value = 4
$ws.Range("Q17").Value = 27
$ws.Range("R17").Value = "'12.9%"
$ws.Range("S17").Value = "'35.0%"

# Row 28 - time changed
$ws.Range("F28").Value = "14:00:00"

# Row 29 - reorder recorder emails
$ws.Range("G29").Value = "marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 33 - reorder recorder emails
$ws.Range("G33").Value = "servinaz@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 40 - reorder recorder emails
$ws.Range("G40").Value = "norhan.mohamed@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"

# Row 60 - reorder recorder emails
$ws.Range("G60").Value = "marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"

# Row 64 - reorder recorder emails
$ws.Range("G64").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 65 - reorder recorder emails
$ws.Range("G65").Value = "majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 71 - reorder recorder emails
$ws.Range("G71").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 89 (PHARMACOLOGY A3 #1) - date changed. Day <= 12 is ambiguous as
# a date, so Excel would silently reinterpret "08/10/2025" as a real
# date serial; force literal text with a leading apostrophe and restore
# the row's own style afterwards (fixed up in the style-restore pass).
$ws.Range("E89").Value = "'08/10/2025"

# Row 90 (PHARMACOLOGY A3 #2) - date/time changed
$ws.Range("E90").Value = "28/10/2025"
$ws.Range("F90").Value = "10:00:00"

# ---------------------------------------------------------------------
# Row 91 (PHYSIOLOGY A3 #1) moved from Pending (yellow) to Recorded
# (green). Copy the visual style from another already-"Recorded" row
# (row 9) onto row 91, then overwrite the date/time/recorder/students/
# status cells with their new values. A/B/C/D keep their text but still
# need the style refresh, which the format-only paste handles.
$ws.Range("A9:I9").Copy() | Out-Null
$ws.Range("A91:I91").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("E91").Value = "29/09/2025"
$ws.Range("F91").Value = "14:00:00"
$ws.Range("G91").Value = "marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("H91").Value = "56/221"
$ws.Range("I91").Value = "Recorded"

# ---------------------------------------------------------------------
# Row 95 - reorder recorder emails
$ws.Range("G95").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 102 - reorder recorder emails
$ws.Range("G102").Value = "Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 121 - time changed
$ws.Range("F121").Value = "14:00:00"

# Row 122 - reorder recorder emails
$ws.Range("G122").Value = "Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# Row 126 - reorder recorder emails
$ws.Range("G126").Value = "servinaz@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 133 - reorder recorder emails
$ws.Range("G133").Value = "Safa.hany@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 137 - date changed (day <= 12, same ambiguous-date issue as E89)
$ws.Range("E137").Value = "'05/11/2025"

# Row 153 - reorder recorder emails
$ws.Range("G153").Value = "Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg"

# Row 157 - reorder recorder emails
$ws.Range("G157").Value = "majorelle.magdy@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"

# Row 158 - reorder recorder emails
$ws.Range("G158").Value = "majorelle.magdy@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg"

# Row 164 - reorder recorder emails
$ws.Range("G164").Value = "Omnia.Mohammed@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, Sara_nabil@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg"

# Row 176 - date/time changed
$ws.Range("E176").Value = "28/10/2025"
$ws.Range("F176").Value = "08:00:00"

# Row 185 - date/time changed
$ws.Range("E185").Value = "13/10/2025"
$ws.Range("F185").Value = "09:00:00"

# ---------------------------------------------------------------------
# Restore the original style index (s="5") on the percentage cells that
# were forced to literal text above; Excel assigns them a fresh
# "quote-prefixed text" style when the apostrophe trick is used, so we
# copy the plain formatting back in from an untouched s="5" donor cell
# (K4), leaving their values untouched.
$ws.Range("K4").Copy() | Out-Null
$ws.Range("L9").PasteSpecial(-4122) | Out-Null
$ws.Range("L10").PasteSpecial(-4122) | Out-Null
$ws.Range("R17").PasteSpecial(-4122) | Out-Null
$ws.Range("S17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Same style restore for the two apostrophe-forced date cells (they
# must stay on the regular "Pending" row style, s="4").
$ws.Range("A89").Copy() | Out-Null
$ws.Range("E89").PasteSpecial(-4122) | Out-Null
$ws.Range("E137").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

Write-Host "Applied daily attendance processing updates."
